$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-07-03 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-07-04 Friday", 2) | Out-Null
$d.Content.Find.Execute("51+4=55", $false, $false, $false, $false, $false, $true, 1, $false, "67-49=18", 2) | Out-Null
$d.Content.Find.Execute("23+57=80", $false, $false, $false, $false, $false, $true, 1, $false, "38+40=78", 2) | Out-Null
$d.Content.Find.Execute("62+10=72", $false, $false, $false, $false, $false, $true, 1, $false, "47+15=62", 2) | Out-Null
$d.Content.Find.Execute("85-12=73", $false, $false, $false, $false, $false, $true, 1, $false, "76+4=80", 2) | Out-Null
$d.Content.Find.Execute("60-41=19", $false, $false, $false, $false, $false, $true, 1, $false, "58+2=60", 2) | Out-Null
$d.Content.Find.Execute("4+70=74", $false, $false, $false, $false, $false, $true, 1, $false, "4+46=50", 2) | Out-Null
$d.Content.Find.Execute("82-20=62", $false, $false, $false, $false, $false, $true, 1, $false, "31-21=10", 2) | Out-Null
$d.Content.Find.Execute("81-11=70", $false, $false, $false, $false, $false, $true, 1, $false, "50+21=71", 2) | Out-Null
$d.Content.Find.Execute("46-21=25", $false, $false, $false, $false, $false, $true, 1, $false, "78-67=11", 2) | Out-Null
$d.Content.Find.Execute("86-63=23", $false, $false, $false, $false, $false, $true, 1, $false, "55+28=83", 2) | Out-Null
$d.Content.Find.Execute("62-20=42", $false, $false, $false, $false, $false, $true, 1, $false, "11+71=82", 2) | Out-Null
$d.Content.Find.Execute("97-28=69", $false, $false, $false, $false, $false, $true, 1, $false, "44-14=30", 2) | Out-Null
$d.Content.Find.Execute("17+75=92", $false, $false, $false, $false, $false, $true, 1, $false, "36+56=92", 2) | Out-Null
$d.Content.Find.Execute("8+8=16", $false, $false, $false, $false, $false, $true, 1, $false, "4+6=10", 2) | Out-Null
$d.Content.Find.Execute("67+13=80", $false, $false, $false, $false, $false, $true, 1, $false, "36+44=80", 2) | Out-Null
$d.Content.Find.Execute("63-56=7", $false, $false, $false, $false, $false, $true, 1, $false, "89-62=27", 2) | Out-Null
$d.Content.Find.Execute("76+16=92", $false, $false, $false, $false, $false, $true, 1, $false, "72-57=15", 2) | Out-Null
$d.Content.Find.Execute("5+88=93", $false, $false, $false, $false, $false, $true, 1, $false, "24+32=56", 2) | Out-Null
$d.Content.Find.Execute("57+17=74", $false, $false, $false, $false, $false, $true, 1, $false, "93-19=74", 2) | Out-Null
$d.Content.Find.Execute("38-20=18", $false, $false, $false, $false, $false, $true, 1, $false, "56-38=18", 2) | Out-Null
$d.Content.Find.Execute("22+71=93", $false, $false, $false, $false, $false, $true, 1, $false, "29+64=93", 2) | Out-Null
$d.Content.Find.Execute("97-59=38", $false, $false, $false, $false, $false, $true, 1, $false, "30+10=40", 2) | Out-Null
$d.Content.Find.Execute("73-1=72", $false, $false, $false, $false, $false, $true, 1, $false, "37+55=92", 2) | Out-Null
$d.Content.Find.Execute("7+55=62", $false, $false, $false, $false, $false, $true, 1, $false, "66+29=95", 2) | Out-Null
$d.Content.Find.Execute("7-1=6", $false, $false, $false, $false, $false, $true, 1, $false, "97-50=47", 2) | Out-Null
$d.Content.Find.Execute("56+25=81", $false, $false, $false, $false, $false, $true, 1, $false, "92-54=38", 2) | Out-Null
$d.Content.Find.Execute("35+42=77", $false, $false, $false, $false, $false, $true, 1, $false, "23+46=69", 2) | Out-Null
$d.Content.Find.Execute("90+3=93", $false, $false, $false, $false, $false, $true, 1, $false, "6+92=98", 2) | Out-Null
$d.Content.Find.Execute("12+67=79", $false, $false, $false, $false, $false, $true, 1, $false, "47+12=59", 2) | Out-Null
$d.Content.Find.Execute("73+17=90", $false, $false, $false, $false, $false, $true, 1, $false, "91-19=72", 2) | Out-Null
$d.Content.Find.Execute("46-40=6", $false, $false, $false, $false, $false, $true, 1, $false, "81-62=19", 2) | Out-Null
$d.Content.Find.Execute("3+84=87", $false, $false, $false, $false, $false, $true, 1, $false, "10-1=9", 2) | Out-Null
$d.Content.Find.Execute("87-72=15", $false, $false, $false, $false, $false, $true, 1, $false, "26+64=90", 2) | Out-Null
$d.Content.Find.Execute("28+44=72", $false, $false, $false, $false, $false, $true, 1, $false, "86-4=82", 2) | Out-Null
$d.Content.Find.Execute("75-2=73", $false, $false, $false, $false, $false, $true, 1, $false, "91+5=96", 2) | Out-Null
$d.Content.Find.Execute("16+21=37", $false, $false, $false, $false, $false, $true, 1, $false, "54-42=12", 2) | Out-Null
$d.Content.Find.Execute("54+26=80", $false, $false, $false, $false, $false, $true, 1, $false, "20+15=35", 2) | Out-Null
$d.Content.Find.Execute("83-69=14", $false, $false, $false, $false, $false, $true, 1, $false, "6+21=27", 2) | Out-Null
$d.Content.Find.Execute("96-36=60", $false, $false, $false, $false, $false, $true, 1, $false, "3+95=98", 2) | Out-Null
$d.Content.Find.Execute("85-34=51", $false, $false, $false, $false, $false, $true, 1, $false, "57-55=2", 2) | Out-Null
$d.Content.Find.Execute("72-69=3", $false, $false, $false, $false, $false, $true, 1, $false, "67-7=60", 2) | Out-Null
$d.Content.Find.Execute("62-57=5", $false, $false, $false, $false, $false, $true, 1, $false, "88+0=88", 2) | Out-Null
$d.Content.Find.Execute("98-44=54", $false, $false, $false, $false, $false, $true, 1, $false, "0+56=56", 2) | Out-Null
$d.Content.Find.Execute("23-4=19", $false, $false, $false, $false, $false, $true, 1, $false, "53-38=15", 2) | Out-Null
$d.Content.Find.Execute("85-85=0", $false, $false, $false, $false, $false, $true, 1, $false, "97-64=33", 2) | Out-Null
$d.Content.Find.Execute("67+0=67", $false, $false, $false, $false, $false, $true, 1, $false, "31-28=3", 2) | Out-Null
$d.Content.Find.Execute("24+12=36", $false, $false, $false, $false, $false, $true, 1, $false, "53+26=79", 2) | Out-Null
$d.Content.Find.Execute("95+4=99", $false, $false, $false, $false, $false, $true, 1, $false, "70-9=61", 2) | Out-Null
$d.Content.Find.Execute("28+61=89", $false, $false, $false, $false, $false, $true, 1, $false, "94+4=98", 2) | Out-Null
$d.Content.Find.Execute("60+27=87", $false, $false, $false, $false, $false, $true, 1, $false, "29-5=24", 2) | Out-Null
$d.Content.Find.Execute("40+14=54", $false, $false, $false, $false, $false, $true, 1, $false, "98-0=98", 2) | Out-Null
$d.Content.Find.Execute("70-37=33", $false, $false, $false, $false, $false, $true, 1, $false, "95-40=55", 2) | Out-Null
$d.Content.Find.Execute("53-42=11", $false, $false, $false, $false, $false, $true, 1, $false, "66-38=28", 2) | Out-Null
$d.Content.Find.Execute("8+53=61", $false, $false, $false, $false, $false, $true, 1, $false, "69-54=15", 2) | Out-Null
$d.Content.Find.Execute("20+17=37", $false, $false, $false, $false, $false, $true, 1, $false, "76-19=57", 2) | Out-Null
$d.Content.Find.Execute("13+56=69", $false, $false, $false, $false, $false, $true, 1, $false, "61+16=77", 2) | Out-Null
$d.Content.Find.Execute("26+12=38", $false, $false, $false, $false, $false, $true, 1, $false, "52-50=2", 2) | Out-Null
$d.Content.Find.Execute("1+21=22", $false, $false, $false, $false, $false, $true, 1, $false, "48-43=5", 2) | Out-Null
$d.Content.Find.Execute("71-22=49", $false, $false, $false, $false, $false, $true, 1, $false, "45-21=24", 2) | Out-Null
$d.Content.Find.Execute("56-37=19", $false, $false, $false, $false, $false, $true, 1, $false, "49-9=40", 2) | Out-Null
$d.Content.Find.Execute("31+55=86", $false, $false, $false, $false, $false, $true, 1, $false, "90-1=89", 2) | Out-Null
$d.Content.Find.Execute("47+24=71", $false, $false, $false, $false, $false, $true, 1, $false, "45+20=65", 2) | Out-Null
$d.Content.Find.Execute("58+1=59", $false, $false, $false, $false, $false, $true, 1, $false, "10+86=96", 2) | Out-Null
$d.Content.Find.Execute("25+36=61", $false, $false, $false, $false, $false, $true, 1, $false, "82+10=92", 2) | Out-Null
$d.Content.Find.Execute("80+16=96", $false, $false, $false, $false, $false, $true, 1, $false, "86+6=92", 2) | Out-Null
$d.Content.Find.Execute("12+43=55", $false, $false, $false, $false, $false, $true, 1, $false, "37+1=38", 2) | Out-Null
$d.Content.Find.Execute("75-66=9", $false, $false, $false, $false, $false, $true, 1, $false, "95-88=7", 2) | Out-Null
$d.Content.Find.Execute("31+0=31", $false, $false, $false, $false, $false, $true, 1, $false, "71-57=14", 2) | Out-Null
$d.Content.Find.Execute("53+0=53", $false, $false, $false, $false, $false, $true, 1, $false, "46+37=83", 2) | Out-Null
$d.Content.Find.Execute("6+93=99", $false, $false, $false, $false, $false, $true, 1, $false, "82-4=78", 2) | Out-Null
$d.Content.Find.Execute("29+51=80", $false, $false, $false, $false, $false, $true, 1, $false, "97-25=72", 2) | Out-Null
$d.Content.Find.Execute("58+21=79", $false, $false, $false, $false, $false, $true, 1, $false, "96-30=66", 2) | Out-Null
$d.Content.Find.Execute("15+1=16", $false, $false, $false, $false, $false, $true, 1, $false, "22+7=29", 2) | Out-Null
$d.Content.Find.Execute("46-0=46", $false, $false, $false, $false, $false, $true, 1, $false, "54+20=74", 2) | Out-Null
$d.Content.Find.Execute("2+69=71", $false, $false, $false, $false, $false, $true, 1, $false, "80-67=13", 2) | Out-Null
$d.Content.Find.Execute("25+56=81", $false, $false, $false, $false, $false, $true, 1, $false, "9-4=5", 2) | Out-Null
$d.Content.Find.Execute("15+54=69", $false, $false, $false, $false, $false, $true, 1, $false, "44+33=77", 2) | Out-Null
$d.Content.Find.Execute("60-30=30", $false, $false, $false, $false, $false, $true, 1, $false, "91-31=60", 2) | Out-Null
$d.Content.Find.Execute("66+2=68", $false, $false, $false, $false, $false, $true, 1, $false, "8+68=76", 2) | Out-Null
$d.Content.Find.Execute("1+3=4", $false, $false, $false, $false, $false, $true, 1, $false, "17+43=60", 2) | Out-Null
$d.Content.Find.Execute("92+7=99", $false, $false, $false, $false, $false, $true, 1, $false, "52-32=20", 2) | Out-Null
$d.Content.Find.Execute("99-24=75", $false, $false, $false, $false, $false, $true, 1, $false, "16+30=46", 2) | Out-Null
$d.Content.Find.Execute("28-19=9", $false, $false, $false, $false, $false, $true, 1, $false, "42+2=44", 2) | Out-Null
$d.Content.Find.Execute("63+17=80", $false, $false, $false, $false, $false, $true, 1, $false, "78-24=54", 2) | Out-Null
$d.Content.Find.Execute("65-38=27", $false, $false, $false, $false, $false, $true, 1, $false, "28+59=87", 2) | Out-Null
$d.Content.Find.Execute("46-16=30", $false, $false, $false, $false, $false, $true, 1, $false, "48-32=16", 2) | Out-Null
$d.Content.Find.Execute("13+55=68", $false, $false, $false, $false, $false, $true, 1, $false, "50-45=5", 2) | Out-Null
$d.Content.Find.Execute("33+8=41", $false, $false, $false, $false, $false, $true, 1, $false, "75+18=93", 2) | Out-Null
$d.Content.Find.Execute("52-14=38", $false, $false, $false, $false, $false, $true, 1, $false, "80-13=67", 2) | Out-Null
$d.Content.Find.Execute("76-7=69", $false, $false, $false, $false, $false, $true, 1, $false, "63+22=85", 2) | Out-Null
$d.Content.Find.Execute("40+18=58", $false, $false, $false, $false, $false, $true, 1, $false, "10+68=78", 2) | Out-Null
$d.Content.Find.Execute("44-15=29", $false, $false, $false, $false, $false, $true, 1, $false, "22+8=30", 2) | Out-Null
$d.Content.Find.Execute("26+30=56", $false, $false, $false, $false, $false, $true, 1, $false, "47-0=47", 2) | Out-Null
$d.Content.Find.Execute("30+69=99", $false, $false, $false, $false, $false, $true, 1, $false, "75-68=7", 2) | Out-Null
$d.Content.Find.Execute("38+29=67", $false, $false, $false, $false, $false, $true, 1, $false, "69-32=37", 2) | Out-Null
$d.Content.Find.Execute("98-92=6", $false, $false, $false, $false, $false, $true, 1, $false, "42-29=13", 2) | Out-Null
$d.Content.Find.Execute("45+49=94", $false, $false, $false, $false, $false, $true, 1, $false, "85-1=84", 2) | Out-Null
$d.Content.Find.Execute("39-23=16", $false, $false, $false, $false, $false, $true, 1, $false, "53+23=76", 2) | Out-Null
$d.Content.Find.Execute("38+5=43", $false, $false, $false, $false, $false, $true, 1, $false, "65+27=92", 2) | Out-Null
$d.Content.Find.Execute("51+39=90", $false, $false, $false, $false, $false, $true, 1, $false, "28+23=51", 2) | Out-Null
